# Weekly refresh of the Espárragos (asparagus) price data for the
# Vega Modelo de Temuco market: rows 2-13 get new Fecha/Variedad/Calidad
# /Volumen/Precio values (the dataset was re-pulled and rows shuffled).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44160
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 1000
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 1000

# Row 3
$ws.Range("D3").Value = 44166
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 285
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1100
$ws.Range("M3").Value = 1054
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 1054

# Row 4
$ws.Range("D4").Value = 44162
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1200
$ws.Range("O4").Value = "Región del Bíobío"
$ws.Range("P4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44162
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 1000
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44165
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 650
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1100
$ws.Range("M6").Value = 1008
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 1008

# Row 7
$ws.Range("D7").Value = 44165
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 800
$ws.Range("M7").Value = 800
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 800

# Row 8
$ws.Range("D8").Value = 44174
$ws.Range("H8").Value = "Verde"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1100
$ws.Range("M8").Value = 1050
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 1050

# Row 9
$ws.Range("D9").Value = 44159
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 2000
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 1000
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1000

# Row 10
$ws.Range("D10").Value = 44168
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 947
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 947

# Row 11
$ws.Range("D11").Value = 44161
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 1000
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1000

# Row 12
$ws.Range("D12").Value = 44175
$ws.Range("H12").Value = "Verde"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1100
$ws.Range("M12").Value = 1067
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1067

# Row 13
$ws.Range("D13").Value = 44167
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 900
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 957
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 957

